# "added error flag for handled exception to TestNG method"
#
# On the "IC" sheet:
#  - Column F ("Execute") is flipped from "no" to "yes" for the test rows
#    (rows 2-71); this marks the handled-exception cases as an error flag
#    that should now execute.
#  - Row 70 gets a "Comment" (column E) of "Bongi" to go with the change.
#  - The view is re-pointed at column F (the column that was just edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IC")
$ws.Activate()

# Flip the "Execute" flag (column F) to "yes" for every data row.
# (Rows 27 and 46 are already "yes" - re-stamping them is a harmless no-op.)
for ($r = 2; $r -le 71; $r++) {
    $ws.Cells.Item($r, 6).Value2 = "yes"
}

# Row 70 previously had an empty Comment cell - fill it in like its
# neighbours (E66:E69), which all read "Bongi".
$ws.Range("E70").Value2 = "Bongi"

# Move the active selection onto the newly-updated column (F2:F71),
# mirroring the frozen-pane scroll position onto column E.
$ws.Range("F2:F71").Select() | Out-Null
